$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-32 and 36-43: update Price (D) and/or Volume(1h) (E) values
$ws.Cells.Item(2, 4).Value = "64.087.43"
$ws.Cells.Item(2, 5).Value = "  -0.28%  "
$ws.Cells.Item(3, 4).Value = "3.472.97"
$ws.Cells.Item(3, 5).Value = "  -0.84%  "
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).Value = "584.91"
$ws.Cells.Item(5, 5).Value = "  -0.28%  "
$ws.Cells.Item(6, 4).Value = "131.59"
$ws.Cells.Item(6, 5).Value = "  -2.12%  "
$ws.Cells.Item(7, 5).Value = "  +0.04%  "
$ws.Cells.Item(8, 4).Value = "0.482"
$ws.Cells.Item(8, 5).Value = "  -0.65%  "
$ws.Cells.Item(9, 4).Value = "7.63"
$ws.Cells.Item(9, 5).Value = "  +4.83%  "
$ws.Cells.Item(10, 5).Value = "  -2.02%  "
$ws.Cells.Item(11, 5).Value = "  +0.10%  "
$ws.Cells.Item(12, 4).Value = "4.070.95"
$ws.Cells.Item(12, 5).Value = "  -0.68%  "
$ws.Cells.Item(13, 5).Value = "  -0.06%  "
$ws.Cells.Item(14, 4).Value = "0.0000177"
$ws.Cells.Item(14, 5).Value = "  -2.57%  "
$ws.Cells.Item(15, 4).Value = "3.476.63"
$ws.Cells.Item(15, 5).Value = "  -0.68%  "
$ws.Cells.Item(16, 4).Value = "64.075.67"
$ws.Cells.Item(16, 5).Value = "  -0.34%  "
$ws.Cells.Item(17, 4).Value = "24.33"
$ws.Cells.Item(17, 5).Value = "  -7.28%  "
$ws.Cells.Item(18, 4).Value = "9.97"
$ws.Cells.Item(18, 5).Value = "  +0.36%  "
$ws.Cells.Item(19, 5).Value = "  -0.41%  "
$ws.Cells.Item(20, 4).Value = "13.45"
$ws.Cells.Item(20, 5).Value = "  -2.19%  "
$ws.Cells.Item(21, 4).Value = "384.31"
$ws.Cells.Item(21, 5).Value = "  -2.37%  "
$ws.Cells.Item(22, 4).Value = "0.574"
$ws.Cells.Item(22, 5).Value = "  +0.23%  "
$ws.Cells.Item(23, 4).Value = "3.615.20"
$ws.Cells.Item(23, 5).Value = "  -0.72%  "
$ws.Cells.Item(24, 4).Value = "74.69"
$ws.Cells.Item(24, 5).Value = "  +0.80%  "
$ws.Cells.Item(25, 5).Value = "  -0.02%  "
$ws.Cells.Item(26, 4).Value = "5.67"
$ws.Cells.Item(26, 5).Value = "  -0.85%  "
$ws.Cells.Item(27, 4).Value = "0.0000112"
$ws.Cells.Item(27, 5).Value = "  -2.33%  "
$ws.Cells.Item(28, 5).Value = "  +0.17%  "
$ws.Cells.Item(29, 5).Value = "  -0.04%  "
$ws.Cells.Item(30, 4).Value = "7.16"
$ws.Cells.Item(30, 5).Value = "  -4.91%  "
$ws.Cells.Item(31, 5).Value = "  -6.20%  "
$ws.Cells.Item(32, 4).Value = "7.93"
$ws.Cells.Item(32, 5).Value = "  -4.40%  "
$ws.Cells.Item(36, 4).Value = "22.93"
$ws.Cells.Item(36, 5).Value = "  -2.40%  "
$ws.Cells.Item(37, 4).Value = "5.17"
$ws.Cells.Item(37, 5).Value = "  -0.78%  "
$ws.Cells.Item(38, 4).Value = "6.77"
$ws.Cells.Item(38, 5).Value = "  -2.36%  "
$ws.Cells.Item(39, 5).Value = "  -3.87%  "
$ws.Cells.Item(40, 4).Value = "162.30"
$ws.Cells.Item(40, 5).Value = "  +0.39%  "
$ws.Cells.Item(41, 4).Value = "0.0775"
$ws.Cells.Item(41, 5).Value = "  -1.12%  "
$ws.Cells.Item(42, 4).Value = "0.798"
$ws.Cells.Item(42, 5).Value = "  -1.04%  "
$ws.Cells.Item(43, 5).Value = "  +0.10%  "

# Rows 33-34, 44-51: coin rows shifted/swapped - update Coin, Link, Price, Volume(1h)
$ws.Cells.Item(33, 2).Value = "RenzoRestakedETH"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Cells.Item(33, 4).Value = "3.504.58"
$ws.Cells.Item(33, 5).Value = "  -0.51%  "
$ws.Cells.Item(34, 2).Value = "Kaspa"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(34, 4).Value = "0.152"
$ws.Cells.Item(34, 5).Value = "  +1.88%  "
$ws.Cells.Item(44, 2).Value = "OKB"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(44, 4).Value = "41.42"
$ws.Cells.Item(44, 5).Value = "  -0.94%  "
$ws.Cells.Item(45, 2).Value = "Filecoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(45, 4).Value = "4.30"
$ws.Cells.Item(45, 5).Value = "  -2.41%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).Value = "23.86"
$ws.Cells.Item(46, 5).Value = "  -5.56%  "
$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).Value = "1.62"
$ws.Cells.Item(47, 5).Value = "  -1.98%  "
$ws.Cells.Item(48, 2).Value = "ONDO"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(48, 4).Value = "1.13"
$ws.Cells.Item(48, 5).Value = "  -3.82%  "
$ws.Cells.Item(49, 2).Value = "SuiNetwork"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(49, 4).Value = "0.913"
$ws.Cells.Item(49, 5).Value = "  +1.92%  "
$ws.Cells.Item(50, 2).Value = "Cosmos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(50, 4).Value = "6.71"
$ws.Cells.Item(50, 5).Value = "  -1.33%  "
$ws.Cells.Item(51, 2).Value = "Maker"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(51, 4).Value = "2.357.09"
$ws.Cells.Item(51, 5).Value = "  -4.69%  "
